$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Project group name: "Axx" -> "A10"
$ws.Range("C2").Value = "A10"

# Header row: D4/E4 were placeholder text "id1"/"id2", now filled with real student ids (numbers)
$ws.Range("D4").Value = 4561619
$ws.Range("E4").Value = 4562771

# Week 3.1 hour entries
# Saturday (row 5): all six students logged 0 hours
$ws.Range("D5:I5").Value = 0

# Sunday (row 6): all six students logged 0 hours
$ws.Range("D6:I6").Value = 0

# Monday (row 7): all six students logged 0 hours
$ws.Range("D7:I7").Value = 0

# Wednesday (row 9): first two students logged hours
$ws.Range("D9").Value = 5
$ws.Range("E9").Value = 4

# Update the selected cell shown in the sheet view
$null = $ws.Range("F10").Select()
